# BCS_template.xlsx - "added the update mode for updating vc numbers"
#
# - "Vendor ID Updates" sheet: rename the 3rd header column from
#   "Vendor ID" to "Updated Vendor ID", bold the header row (new style),
#   resize columns B/C, add a portrait page setup, and move the selection.
# - "New Vendors" sheet: just move the remembered selection back to A2.

$wb = $excel.ActiveWorkbook

$wsNew     = $wb.Worksheets.Item("New Vendors")
$wsUpdates = $wb.Worksheets.Item("Vendor ID Updates")

# --- "Vendor ID Updates" sheet ------------------------------------------

# Rename the header so it's clear this column now holds the *updated*
# vendor/V-Commerce id once it has been refreshed.
$wsUpdates.Range("C1").Value = "Updated Vendor ID"

# Make the whole header row bold (adds a new bold, non-centered cell style).
$wsUpdates.Range("A1:C1").Font.Bold = $true

# Widen/narrow the Vendor Name / Updated Vendor ID columns to fit the new text.
$wsUpdates.Columns.Item(2).ColumnWidth = 16.25
$wsUpdates.Columns.Item(3).ColumnWidth = 17.875

# Give the sheet an explicit (portrait) print setup.
$wsUpdates.PageSetup.Orientation = 1

# Remember a new cursor position on this sheet (it is not the active tab).
$wsUpdates.Range("C7").Select() | Out-Null

# --- "New Vendors" sheet -------------------------------------------------

# This is the active tab, so select it last - re-point the remembered
# selection at A2.
$wsNew.Activate() | Out-Null
$wsNew.Range("A2").Select() | Out-Null
